# Updated cryptos list on Fri Mar 15 07:35:12 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the "Price" cells that look like plain decimal numbers to be stored
# as text (matching the source data, which keeps prices as literal strings,
# e.g. "581.34"), so Excel does not silently coerce them into doubles.
$textRows = @(5, 6, 10, 12, 14, 18, 19, 20, 22, 23, 27, 28, 30, 31, 32, 33, 35, 36, 37, 39, 42, 45, 46, 47, 49, 51)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.267.31"
$ws.Range("E2").Value = "  -6.85%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.737.54"
$ws.Range("E3").Value = "  -5.73%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "581.34"
$ws.Range("E5").Value = "  -4.93%  "

# Row 6 - Solana
$ws.Range("D6").Value = "175.13"
$ws.Range("E6").Value = "  +4.02%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.730.76"
$ws.Range("E7").Value = "  -5.69%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -6.72%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.07%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "0.716"
$ws.Range("E10").Value = "  -6.13%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -9.22%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "53.74"
$ws.Range("E12").Value = "  -3.76%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -10.07%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "10.78"
$ws.Range("E14").Value = "  -3.68%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.331.27"
$ws.Range("E15").Value = "  -6.04%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.762.46"
$ws.Range("E16").Value = "  -5.34%  "

# Row 17 - Chainlink
$ws.Range("E17").Value = "  -4.30%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "13.17"
$ws.Range("E18").Value = "  -6.67%  "

# Row 19 - was Polygon, becomes TRON
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "0.127"
$ws.Range("E19").Value = "  -2.75%  "

# Row 20 - was TRON, becomes Polygon
$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D20").Value = "1.15"
$ws.Range("E20").Value = "  -7.01%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "68.165.48"
$ws.Range("E21").Value = "  -6.85%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "412.84"
$ws.Range("E22").Value = "  -5.70%  "

# Row 23 - PancakeSwap
$ws.Range("D23").Value = "4.55"
$ws.Range("E23").Value = "  -5.64%  "

# Row 24 - Litecoin
$ws.Range("E24").Value = "  -6.66%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -7.65%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("E26").Value = "  -8.42%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "10.84"
$ws.Range("E27").Value = "  -1.67%  "

# Row 28 - Toncoin
$ws.Range("D28").Value = "3.86"
$ws.Range("E28").Value = "  -5.23%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  +0.80%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "9.64"
$ws.Range("E30").Value = "  -7.56%  "

# Row 31 - was NEARProtocol, becomes EthereumClassic
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "33.23"
$ws.Range("E31").Value = "  -7.59%  "

# Row 32 - was EthereumClassic, becomes NEARProtocol
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "7.97"
$ws.Range("E32").Value = "  +1.89%  "

# Row 33 - Cosmos
$ws.Range("D33").Value = "12.80"
$ws.Range("E33").Value = "  -7.78%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -7.96%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "615.39"
$ws.Range("E35").Value = "  -4.88%  "

# Row 36 - OKB
$ws.Range("D36").Value = "65.60"
$ws.Range("E36").Value = "  -6.62%  "

# Row 37 - InjectiveProtocol
$ws.Range("D37").Value = "44.01"
$ws.Range("E37").Value = "  -7.88%  "

# Row 38 - PEPE
$ws.Range("E38").Value = "  -11.46%  "

# Row 39 - TheGraph
$ws.Range("D39").Value = "0.405"
$ws.Range("E39").Value = "  -5.26%  "

# Row 40 - Dai
$ws.Range("E40").Value = "  +0.32%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  -0.15%  "

# Row 42 - dogwifhat
$ws.Range("D42").Value = "3.21"
$ws.Range("E42").Value = "  +2.31%  "

# Row 43 - Kaspa
$ws.Range("E43").Value = "  -5.50%  "

# Row 44 - ThetaToken
$ws.Range("E44").Value = "  -8.47%  "

# Row 45 - VeChain
$ws.Range("D45").Value = "0.0449"
$ws.Range("E45").Value = "  -7.03%  "

# Row 46 - Fetch.AI
$ws.Range("D46").Value = "2.63"
$ws.Range("E46").Value = "  +2.74%  "

# Row 47 - THORChain
$ws.Range("D47").Value = "9.44"
$ws.Range("E47").Value = "  -10.02%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -7.90%  "

# Row 49 - WEMIXToken
$ws.Range("D49").Value = "2.73"
$ws.Range("E49").Value = "  -14.75%  "

# Row 50 - was ApeXProtocol, becomes Maker
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.745.30"
$ws.Range("E50").Value = "  -1.90%  "

# Row 51 - was Maker, becomes ApeXProtocol
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "3.12"
$ws.Range("E51").Value = "  -9.11%  "
